# Apply updated Pseudo-Voigt fit results for rows 2-10 (Ne calibration tweak)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9978598860796123
$ws.Range("D2").Value = 0.00003987360630922584
$ws.Range("E2").Value = 331.1864109950466
$ws.Range("F2").Value = 1447.507791711392
$ws.Range("G2").Value = 1116.321380716346
$ws.Range("H2").Value = 49150.03938907994
$ws.Range("I2").Value = 4064.558187157941
$ws.Range("J2").Value = 122.7893778732266
$ws.Range("K2").Value = 11.05666432640984
$ws.Range("L2").Value = 0.00003987360630922584
$ws.Range("M2").Value = 0.3831035754123102
$ws.Range("O2").Value = 0.00257
$ws.Range("P2").Value = 0.405291536749454
$ws.Range("Q2").Value = 0.4051926376392343
$ws.Range("S2").Value = 0.01295310312974792
$ws.Range("T2").Value = 0.2777918959532876
$ws.Range("U2").Value = 0.9978131173855738
$ws.Range("V2").Value = 0.9979066591580656
$ws.Range("W2").Value = 133.8460421996364

$ws.Range("C3").Value = 0.9978264842220211
$ws.Range("D3").Value = 0.0000357430383875218
$ws.Range("E3").Value = 331.1974973327126
$ws.Range("F3").Value = 1447.503005381121
$ws.Range("G3").Value = 1116.305508048408
$ws.Range("H3").Value = 48194.63294913297
$ws.Range("I3").Value = 4249.843322072755
$ws.Range("J3").Value = 113.0786756227321
$ws.Range("K3").Value = 9.905378667281072
$ws.Range("L3").Value = 0.0000357430383875218
$ws.Range("M3").Value = 0.3857421092605332
$ws.Range("P3").Value = 0.3997918883762511
$ws.Range("Q3").Value = 0.4092017687249572
$ws.Range("S3").Value = 0.01159216800643426
$ws.Range("T3").Value = 0.2893279115707663
$ws.Range("U3").Value = 0.9977843306245427
$ws.Range("V3").Value = 0.9978686413813932
$ws.Range("W3").Value = 122.9840542900132

$ws.Range("C4").Value = 0.9978100509622898
$ws.Range("D4").Value = 0.00003933335809059285
$ws.Range("E4").Value = 331.2029519325042
$ws.Range("F4").Value = 1447.521962191433
$ws.Range("G4").Value = 1116.319010258929
$ws.Range("H4").Value = 46080.58188972915
$ws.Range("I4").Value = 3842.899468127829
$ws.Range("J4").Value = 109.3326919313476
$ws.Range("K4").Value = 10.20425453009774
$ws.Range("L4").Value = 0.00003933335809059284
$ws.Range("M4").Value = 0.3896353958047211
$ws.Range("P4").Value = 0.4037278540818179
$ws.Range("Q4").Value = 0.406028126976008
$ws.Range("S4").Value = 0.01279868269207621
$ws.Range("T4").Value = 0.2818591713753829
$ws.Range("U4").Value = 0.9977641738530534
$ws.Range("V4").Value = 0.997855932290571
$ws.Range("W4").Value = 119.5369464614454

$ws.Range("C5").Value = 0.9977826198315888
$ws.Range("D5").Value = 0.00004053693854759805
$ws.Range("E5").Value = 331.2120573942375
$ws.Range("F5").Value = 1447.539256518417
$ws.Range("G5").Value = 1116.32719912418
$ws.Range("H5").Value = 44440.9420876005
$ws.Range("I5").Value = 3811.741538002881
$ws.Range("J5").Value = 106.9637865077495
$ws.Range("K5").Value = 10.36403580892033
$ws.Range("L5").Value = 0.00004053693854759805
$ws.Range("M5").Value = 0.3914122320374968
$ws.Range("P5").Value = 0.4092346608981216
$ws.Range("Q5").Value = 0.4082686204520876
$ws.Range("S5").Value = 0.0132064281462208
$ws.Range("T5").Value = 0.2859551890075377
$ws.Range("U5").Value = 0.997735547140579
$ws.Range("V5").Value = 0.9978296969645427
$ws.Range("W5").Value = 117.3278223166699

$ws.Range("C6").Value = 0.9977781961827885
$ws.Range("D6").Value = 0.00003799534336316177
$ws.Range("E6").Value = 331.2135258226176
$ws.Range("F6").Value = 1447.541011877257
$ws.Range("G6").Value = 1116.32748605464
$ws.Range("H6").Value = 45288.66855610132
$ws.Range("I6").Value = 3947.474007742397
$ws.Range("J6").Value = 110.8975282102296
$ws.Range("K6").Value = 10.91844694552933
$ws.Range("L6").Value = 0.00003799534336316177
$ws.Range("M6").Value = 0.3933158774344163
$ws.Range("P6").Value = 0.4163762948292194
$ws.Range("Q6").Value = 0.4098356067117886
$ws.Range("S6").Value = 0.01233172102207001
$ws.Range("T6").Value = 0.2926682750136492
$ws.Range("U6").Value = 0.997733487615888
$ws.Range("V6").Value = 0.9978229087566619
$ws.Range("W6").Value = 121.8159751557589

$ws.Range("C7").Value = 0.997775737608652
$ws.Range("D7").Value = 0.00003521989263082476
$ws.Range("E7").Value = 331.2143419509093
$ws.Range("F7").Value = 1447.539632678197
$ws.Range("G7").Value = 1116.325290727288
$ws.Range("H7").Value = 46370.34161343362
$ws.Range("I7").Value = 4254.113932483191
$ws.Range("J7").Value = 116.1392139824934
$ws.Range("K7").Value = 11.1446659173381
$ws.Range("L7").Value = 0.00003521989263082476
$ws.Range("M7").Value = 0.3925913687643211
$ws.Range("P7").Value = 0.4170039054175857
$ws.Range("Q7").Value = 0.4046799468193902
$ws.Range("S7").Value = 0.01137871289278358
$ws.Range("T7").Value = 0.3126667295924849
$ws.Range("U7").Value = 0.9977337191951056
$ws.Range("V7").Value = 0.9978177595614622
$ws.Range("W7").Value = 127.2838798998315

$ws.Range("C8").Value = 0.9977600522180611
$ws.Range("D8").Value = 0.0000396303108442229
$ws.Range("E8").Value = 331.2195488403927
$ws.Range("F8").Value = 1447.534557918229
$ws.Range("G8").Value = 1116.315009077836
$ws.Range("H8").Value = 44105.49495871713
$ws.Range("I8").Value = 3915.579106580432
$ws.Range("J8").Value = 102.6329463068532
$ws.Range("K8").Value = 11.22393506653637
$ws.Range("L8").Value = 0.0000396303108442229
$ws.Range("M8").Value = 0.3949879254397151
$ws.Range("P8").Value = 0.4058043278678481
$ws.Range("Q8").Value = 0.4070156761523468
$ws.Range("S8").Value = 0.01291243725375591
$ws.Range("T8").Value = 0.3022478161670795
$ws.Range("U8").Value = 0.9977140479129187
$ws.Range("V8").Value = 0.9978060607658893
$ws.Range("W8").Value = 113.8568813733896

$ws.Range("C9").Value = 0.9977557132171323
$ws.Range("D9").Value = 0.0000408649213211015
$ws.Range("E9").Value = 331.2209892349813
$ws.Range("F9").Value = 1447.516894431486
$ws.Range("G9").Value = 1116.295905196504
$ws.Range("H9").Value = 43261.39606674626
$ws.Range("I9").Value = 3835.072595910358
$ws.Range("J9").Value = 108.0927717612777
$ws.Range("K9").Value = 11.09514505542197
$ws.Range("L9").Value = 0.0000408649213211015
$ws.Range("M9").Value = 0.3929697387808552
$ws.Range("P9").Value = 0.4140404255556192
$ws.Range("Q9").Value = 0.4073028395439339
$ws.Range("S9").Value = 0.01328909246086156
$ws.Range("T9").Value = 0.3142053281481823
$ws.Range("U9").Value = 0.9977079422635405
$ws.Range("V9").Value = 0.9978034887455564
$ws.Range("W9").Value = 119.1879168166997

$ws.Range("C10").Value = 0.9977568843223
$ws.Range("D10").Value = 0.00004065606820687621
$ws.Range("E10").Value = 331.2206004683205
$ws.Range("F10").Value = 1447.524736945334
$ws.Range("G10").Value = 1116.304136477014
$ws.Range("H10").Value = 43156.99290066765
$ws.Range("I10").Value = 3790.648654420056
$ws.Range("J10").Value = 103.7940358126246
$ws.Range("K10").Value = 11.52920847606356
$ws.Range("L10").Value = 0.00004065606820687622
$ws.Range("M10").Value = 0.3925621523975833
$ws.Range("P10").Value = 0.4168855850620845
$ws.Range("Q10").Value = 0.4084576872680725
$ws.Range("S10").Value = 0.01324322412067986
$ws.Range("T10").Value = 0.3084446521308241
$ws.Range("U10").Value = 0.9977096429855106
$ws.Range("V10").Value = 0.9978041301330356
$ws.Range("W10").Value = 115.3232442886881

Write-Output "Updated 163 cells"
